$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.033299267841694
$ws.Range("D2").Value2 = 1.040530985516091
$ws.Range("E2").Value2 = 1.032632998719386
$ws.Range("F2").Value2 = 1.047805281267926
$ws.Range("I2").Value2 = 1.033253539088952
$ws.Range("J2").Value2 = 1.038424756301964
$ws.Range("K2").Value2 = 1.043313070620796
$ws.Range("L2").Value2 = 1.035437667400937
$ws.Range("M2").Value2 = 1.050566892918201
$ws.Range("N2").Value2 = 1.039899438564897

$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.034311450429669
$ws.Range("D3").Value2 = 1.041462593743852
$ws.Range("E3").Value2 = 1.033493724896626
$ws.Range("F3").Value2 = 1.048897026253906
$ws.Range("I3").Value2 = 1.03340583016553
$ws.Range("J3").Value2 = 1.039079276549314
$ws.Range("K3").Value2 = 1.044054802476562
$ws.Range("L3").Value2 = 1.036107067550997
$ws.Range("M3").Value2 = 1.051469831016945
$ws.Range("N3").Value2 = 1.040554888306074

$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.034966143543362
$ws.Range("D4").Value2 = 1.042065524778049
$ws.Range("E4").Value2 = 1.034050844024983
$ws.Range("F4").Value2 = 1.049603890632234
$ws.Range("I4").Value2 = 1.033502094869176
$ws.Range("J4").Value2 = 1.039501958895028
$ws.Range("K4").Value2 = 1.044534246544083
$ws.Range("L4").Value2 = 1.036539747952801
$ws.Range("M4").Value2 = 1.052053932386183
$ws.Range("N4").Value2 = 1.0409781709092

$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.035241314851391
$ws.Range("D5").Value2 = 1.042319024625562
$ws.Range("E5").Value2 = 1.034285097228027
$ws.Range("F5").Value2 = 1.049901159462006
$ws.Range("I5").Value2 = 1.033542018624184
$ws.Range("J5").Value2 = 1.039679453852148
$ws.Range("K5").Value2 = 1.044735683132684
$ws.Range("L5").Value2 = 1.036721534568381
$ws.Range("M5").Value2 = 1.052299449803559
$ws.Range("N5").Value2 = 1.041155917929518

$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.03528751370148
$ws.Range("D6").Value2 = 1.04236158996692
$ws.Range("E6").Value2 = 1.034324431712844
$ws.Range("F6").Value2 = 1.049951078235541
$ws.Range("I6").Value2 = 1.033548689966848
$ws.Range("J6").Value2 = 1.03970924425092
$ws.Range("K6").Value2 = 1.044769498069728
$ws.Range("L6").Value2 = 1.036752050742412
$ws.Range("M6").Value2 = 1.052340670968321
$ws.Range("N6").Value2 = 1.041185750634073

$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.03496982063849
$ws.Range("D7").Value2 = 1.042068911947658
$ws.Range("E7").Value2 = 1.034053973970325
$ws.Range("F7").Value2 = 1.049607862348588
$ws.Range("I7").Value2 = 1.033502630478255
$ws.Range("J7").Value2 = 1.03950433138067
$ws.Range("K7").Value2 = 1.044536938628537
$ws.Range("L7").Value2 = 1.036542177436201
$ws.Range("M7").Value2 = 1.052057213155031
$ws.Range("N7").Value2 = 1.040980546764044

$ws.Range("B8").Value2 = 1.019999999999999
$ws.Range("C8").Value2 = 1.033641392707932
$ws.Range("D8").Value2 = 1.04084580200691
$ws.Range("E8").Value2 = 1.032923849128113
$ws.Range("F8").Value2 = 1.048174151795515
$ws.Range("I8").Value2 = 1.033305477839786
$ws.Range("J8").Value2 = 1.038646127256515
$ws.Range("K8").Value2 = 1.043563846825609
$ws.Range("L8").Value2 = 1.035663990497876
$ws.Range("M8").Value2 = 1.050872077896477
$ws.Range("N8").Value2 = 1.040121123891595

$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.031298575745147
$ws.Range("D9").Value2 = 1.038691459755404
$ws.Range("E9").Value2 = 1.030933774313417
$ws.Range("F9").Value2 = 1.045651101718779
$ws.Range("I9").Value2 = 1.032940650903305
$ws.Range("J9").Value2 = 1.037127482764415
$ws.Range("K9").Value2 = 1.041845282570388
$ws.Range("L9").Value2 = 1.034112967826453
$ws.Range("M9").Value2 = 1.048782514473169
$ws.Range("N9").Value2 = 1.038600322750125

$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.029735390046008
$ws.Range("D10").Value2 = 1.037255894759414
$ws.Range("E10").Value2 = 1.029608008868507
$ws.Range("F10").Value2 = 1.043971331141264
$ws.Range("I10").Value2 = 1.032685755695967
$ws.Range("J10").Value2 = 1.03611079359996
$ws.Range("K10").Value2 = 1.040697012162863
$ws.Range("L10").Value2 = 1.033076601812777
$ws.Range("M10").Value2 = 1.047388688779178
$ws.Range("N10").Value2 = 1.037582189770441

$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.029058203754013
$ws.Range("D11").Value2 = 1.036634443128664
$ws.Range("E11").Value2 = 1.029034171838503
$ws.Range("F11").Value2 = 1.04324451246632
$ws.Range("I11").Value2 = 1.032572620282212
$ws.Range("J11").Value2 = 1.035669551088786
$ws.Range("K11").Value2 = 1.040199196290128
$ws.Range("L11").Value2 = 1.03262729243882
$ws.Range("M11").Value2 = 1.046784965554074
$ws.Range("N11").Value2 = 1.03714032064429

$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.028806618960333
$ws.Range("D12").Value2 = 1.036403632414595
$ws.Range("E12").Value2 = 1.028821058092924
$ws.Range("F12").Value2 = 1.042974619928853
$ws.Range("I12").Value2 = 1.032530181944824
$ws.Range("J12").Value2 = 1.035505502664847
$ws.Range("K12").Value2 = 1.04001419441449
$ws.Range("L12").Value2 = 1.032460315803722
$ws.Range("M12").Value2 = 1.046560687952294
$ws.Range("N12").Value2 = 1.036976039252772

$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.028860586910567
$ws.Range("D13").Value2 = 1.036453140999236
$ws.Range("E13").Value2 = 1.028866770126975
$ws.Range("F13").Value2 = 1.043032509158466
$ws.Range("I13").Value2 = 1.032539303875744
$ws.Range("J13").Value2 = 1.035540698450327
$ws.Range("K13").Value2 = 1.040053882052794
$ws.Range("L13").Value2 = 1.032496136613957
$ws.Range("M13").Value2 = 1.04660879751842
$ws.Range("N13").Value2 = 1.037011285020305

$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.029037408649153
$ws.Range("D14").Value2 = 1.03661536374911
$ws.Range("E14").Value2 = 1.029016555074888
$ws.Range("F14").Value2 = 1.043222201426205
$ws.Range("I14").Value2 = 1.032569120775275
$ws.Range("J14").Value2 = 1.035655993879605
$ws.Range("K14").Value2 = 1.040183905829727
$ws.Range("L14").Value2 = 1.032613491791508
$ws.Range("M14").Value2 = 1.046766427265562
$ws.Range("N14").Value2 = 1.037126744182317

$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.029146348040769
$ws.Range("D15").Value2 = 1.036715317744473
$ws.Range("E15").Value2 = 1.029108847164206
$ws.Range("F15").Value2 = 1.043339087736503
$ws.Range("I15").Value2 = 1.032587436998999
$ws.Range("J15").Value2 = 1.035727011158954
$ws.Range("K15").Value2 = 1.040264005735278
$ws.Range("L15").Value2 = 1.032685787186515
$ws.Range("M15").Value2 = 1.046863544463669
$ws.Range("N15").Value2 = 1.037197862314348

$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.029780325992505
$ws.Range("D16").Value2 = 1.037297141812504
$ws.Range("E16").Value2 = 1.029646097400679
$ws.Range("F16").Value2 = 1.044019578950756
$ws.Range("I16").Value2 = 1.032693205924434
$ws.Range("J16").Value2 = 1.036140056174781
$ws.Range("K16").Value2 = 1.040730037798075
$ws.Range("L16").Value2 = 1.033106409316189
$ws.Range("M16").Value2 = 1.047428751972802
$ws.Range("N16").Value2 = 1.037611493901475

$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.030177918795251
$ws.Range("D17").Value2 = 1.037662147149668
$ws.Range("E17").Value2 = 1.029983161790307
$ws.Range("F17").Value2 = 1.044446575967463
$ws.Range("I17").Value2 = 1.032758812281044
$ws.Range("J17").Value2 = 1.036398878401162
$ws.Range("K17").Value2 = 1.041022205068177
$ws.Range("L17").Value2 = 1.033370105888663
$ws.Range("M17").Value2 = 1.047783241693826
$ws.Range("N17").Value2 = 1.037870683685105

$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.030409797412701
$ws.Range("D18").Value2 = 1.037875063754135
$ws.Range("E18").Value2 = 1.030179787720211
$ws.Range("F18").Value2 = 1.044695687524172
$ws.Range("I18").Value2 = 1.032796812631806
$ws.Range("J18").Value2 = 1.036549747612678
$ws.Range("K18").Value2 = 1.041192562685199
$ws.Range("L18").Value2 = 1.033523861959973
$ws.Range("M18").Value2 = 1.047989991442798
$ws.Range("N18").Value2 = 1.038021767148206

$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.030488856851939
$ws.Range("D19").Value2 = 1.037947665373276
$ws.Range("E19").Value2 = 1.030246835789847
$ws.Range("F19").Value2 = 1.044780636843327
$ws.Range("I19").Value2 = 1.032809724513455
$ws.Range("J19").Value2 = 1.036601173581339
$ws.Range("K19").Value2 = 1.04125064027894
$ws.Range("L19").Value2 = 1.03357627968948
$ws.Range("M19").Value2 = 1.048060484685027
$ws.Range("N19").Value2 = 1.038073266147641

$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.030135264009451
$ws.Range("D20").Value2 = 1.037622983959362
$ws.Range("E20").Value2 = 1.02994699569392
$ws.Range("F20").Value2 = 1.044400757908368
$ws.Range("I20").Value2 = 1.032751800925247
$ws.Range("J20").Value2 = 1.036371119299169
$ws.Range("K20").Value2 = 1.040990864344785
$ws.Range("L20").Value2 = 1.033341819291122
$ws.Range("M20").Value2 = 1.047745210172617
$ws.Range("N20").Value2 = 1.037842885162003

$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.028985340370541
$ws.Range("D21").Value2 = 1.036567592524097
$ws.Range("E21").Value2 = 1.02897244617346
$ws.Range("F21").Value2 = 1.043166339563683
$ws.Range("I21").Value2 = 1.032560351888531
$ws.Range("J21").Value2 = 1.035622046423519
$ws.Range("K21").Value2 = 1.040145619574714
$ws.Range("L21").Value2 = 1.032578935908192
$ws.Range("M21").Value2 = 1.046720010002435
$ws.Range("N21").Value2 = 1.03709274851695

$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.028262061981331
$ws.Range("D22").Value2 = 1.035904165634923
$ws.Range("E22").Value2 = 1.028359910290707
$ws.Range("F22").Value2 = 1.042390676321714
$ws.Range("I22").Value2 = 1.032437580557832
$ws.Range("J22").Value2 = 1.035150199311321
$ws.Range("K22").Value2 = 1.039613654906067
$ws.Range("L22").Value2 = 1.032098800091107
$ws.Range("M22").Value2 = 1.046075264577106
$ws.Range("N22").Value2 = 1.036620231327731

$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.02864551171172
$ws.Range("D23").Value2 = 1.036255847407717
$ws.Range("E23").Value2 = 1.028684607799043
$ws.Range("F23").Value2 = 1.042801825984713
$ws.Range("I23").Value2 = 1.032502891275341
$ws.Range("J23").Value2 = 1.035400417246662
$ws.Range("K23").Value2 = 1.039895709167258
$ws.Range("L23").Value2 = 1.032353374606844
$ws.Range("M23").Value2 = 1.046417071540862
$ws.Range("N23").Value2 = 1.03687080460124

$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.030154537965271
$ws.Range("D24").Value2 = 1.037640680072363
$ws.Range("E24").Value2 = 1.029963337528195
$ws.Range("F24").Value2 = 1.044421460958044
$ws.Range("I24").Value2 = 1.032754969879263
$ws.Range("J24").Value2 = 1.0363836627434
$ws.Range("K24").Value2 = 1.041005026049649
$ws.Range("L24").Value2 = 1.033354600952769
$ws.Range("M24").Value2 = 1.047762395035844
$ws.Range("N24").Value2 = 1.037855446419363

$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.031904481149073
$ws.Range("D25").Value2 = 1.039248294531522
$ws.Range("E25").Value2 = 1.031448091882064
$ws.Range("F25").Value2 = 1.046302973468874
$ws.Range("I25").Value2 = 1.033037028066878
$ws.Range("J25").Value2 = 1.037520841654797
$ws.Range("K25").Value2 = 1.042290026093973
$ws.Range("L25").Value2 = 1.034514361079722
$ws.Range("M25").Value2 = 1.049322856645991
$ws.Range("N25").Value2 = 1.038994240255251
